# Weekly update: insert a new price record for "Albahaca" (row 63) and
# push the existing historical rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63; this shifts rows 63-77 down to 64-78
# and automatically extends the used range (dimension) to row 78.
$ws.Rows(63).Insert()

# Populate the newly inserted row 63 with the latest weekly record.
$ws.Range("A63").Value = 1
$ws.Range("B63").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C63").Value = "Arica y Parinacota"
$ws.Range("D63").Value = 45211
$ws.Range("D63").NumberFormat = $ws.Range("D64").NumberFormat
$ws.Range("E63").Value = 15
$ws.Range("F63").Value = 100112052
$ws.Range("G63").Value = "Albahaca"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 550
$ws.Range("K63").Value = 1400
$ws.Range("L63").Value = 1500
$ws.Range("M63").Value = 1455
$ws.Range("N63").Value = "$/paquete"
$ws.Range("O63").Value = "Región de Arica y Parinacota"
$ws.Range("P63").Value = 1455
$ws.Range("Q63").Value = 1
$ws.Range("R63").Value = "Hortaliza"
